# Applies the "Update countries & provincias Spain" COVID dataset refresh:
#   1. Fixes the display order of a handful of country names that were
#      shifted in the shared-string table (Suazilandia/Yemen/Nueva Zelanda,
#      Benin/Mozambique, Seychelles/Aruba/Barbados).
#   2. Refreshes the case/death/recovery counters for the affected rows.
#   3. Bumps the "Datos actualizados..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---- Country name corrections (column A) ----
$ws.Range("A128").Value = "Suazilandia"
$ws.Range("A129").Value = "Yemen"
$ws.Range("A130").Value = "Nueva Zelanda"
$ws.Range("A132").Value = "Benin"
$ws.Range("A133").Value = "Mozambique"
$ws.Range("A183").Value = "Seychelles"
$ws.Range("A184").Value = "Aruba"
$ws.Range("A185").Value = "Barbados"

# ---- Updated statistics (columns B-H) ----
$ws.Range("B4").Value = 3683324
$ws.Range("C4").Value = 66497
$ws.Range("D4").Value = 1673175
$ws.Range("E4").Value = 1869133
$ws.Range("G4").Value = 872
$ws.Range("H4").Value = 141016
$ws.Range("B5").Value = 2014738
$ws.Range("C5").Value = 43829
$ws.Range("E5").Value = 571141
$ws.Range("G5").Value = 1299
$ws.Range("H5").Value = 76822
$ws.Range("B19").Value = 201836
$ws.Range("C19").Value = 584
$ws.Range("E19").Value = 6279
$ws.Range("B24").Value = 109264
$ws.Range("C24").Value = 435
$ws.Range("D24").Value = 72836
$ws.Range("E24").Value = 27603
$ws.Range("B51").Value = 34854
$ws.Range("C51").Value = 595
$ws.Range("D51").Value = 14292
$ws.Range("E51").Value = 19793
$ws.Range("G51").Value = 9
$ws.Range("H51").Value = 769
$ws.Range("B54").Value = 32939
$ws.Range("C54").Value = 865
$ws.Range("D54").Value = 4807
$ws.Range("E54").Value = 26728
$ws.Range("G54").Value = 54
$ws.Range("H54").Value = 1404
$ws.Range("B92").Value = 6359
$ws.Range("C92").Value = 83
$ws.Range("D92").Value = 5012
$ws.Range("E92").Value = 1308
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 39
$ws.Range("B95").Value = 5659
$ws.Range("C95").Value = 95
$ws.Range("D95").Value = 2993
$ws.Range("E95").Value = 2516
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 150
$ws.Range("B98").Value = 4373
$ws.Range("C98").Value = 11
$ws.Range("D98").Value = 1265
$ws.Range("E98").Value = 3055
$ws.Range("B106").Value = 3106
$ws.Range("C106").Value = 23
$ws.Range("D106").Value = 1444
$ws.Range("E106").Value = 1569
$ws.Range("B108").Value = 2899
$ws.Range("C108").Value = 68
$ws.Range("D108").Value = 2339
$ws.Range("E108").Value = 545
$ws.Range("B109").Value = 2778
$ws.Range("C109").Value = 35
$ws.Range("E109").Value = 160
$ws.Range("B115").Value = 2358
$ws.Range("C115").Value = 136
$ws.Range("E115").Value = 1721
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 48
$ws.Range("B116").Value = 2171
$ws.Range("C116").Value = 18
$ws.Range("E116").Value = 955
$ws.Range("B125").Value = 1678
$ws.Range("C125").Value = 10
$ws.Range("D125").Value = 1213
$ws.Range("E125").Value = 401
$ws.Range("C128").Value = 63
$ws.Range("D128").Value = 736
$ws.Range("E128").Value = 795
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 21
$ws.Range("B129").Value = 1552
$ws.Range("C129").Value = 26
$ws.Range("D129").Value = 695
$ws.Range("E129").Value = 419
$ws.Range("G129").Value = 5
$ws.Range("H129").Value = 438
$ws.Range("B130").Value = 1548
$ws.Range("C130").Value = 1
$ws.Range("D130").Value = 1499
$ws.Range("E130").Value = 27
$ws.Range("H130").Value = 22
$ws.Range("B132").Value = 1463
$ws.Range("C132").Value = 85
$ws.Range("D132").Value = 557
$ws.Range("E132").Value = 878
$ws.Range("G132").Value = 2
$ws.Range("H132").Value = 28
$ws.Range("B133").Value = 1383
$ws.Range("C133").Value = 53
$ws.Range("D133").Value = 375
$ws.Range("E133").Value = 999
$ws.Range("H133").Value = 9
$ws.Range("B140").Value = 1070
$ws.Range("C140").Value = 14
$ws.Range("D140").Value = 486
$ws.Range("E140").Value = 516
$ws.Range("G140").Value = 17
$ws.Range("H140").Value = 68
$ws.Range("B144").Value = 1031
$ws.Range("C144").Value = 6
$ws.Range("D144").Value = 845
$ws.Range("B145").Value = 1026
$ws.Range("C145").Value = 17
$ws.Range("D145").Value = 916
$ws.Range("E145").Value = 78
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 32
$ws.Range("B147").Value = 904
$ws.Range("C147").Value = 67
$ws.Range("E147").Value = 305
$ws.Range("B149").Value = 877
$ws.Range("C149").Value = 15
$ws.Range("E149").Value = 22
$ws.Range("B152").Value = 740
$ws.Range("C152").Value = 3
$ws.Range("D152").Value = 325
$ws.Range("B164").Value = 339
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 270
$ws.Range("E164").Value = 63
$ws.Range("B172").Value = 251
$ws.Range("C172").Value = 19
$ws.Range("D172").Value = 149
$ws.Range("E172").Value = 102
$ws.Range("B183").Value = 108
$ws.Range("C183").Value = 8
$ws.Range("D183").Value = 27
$ws.Range("E183").Value = 81
$ws.Range("H183").Value = 0
$ws.Range("B184").Value = 106
$ws.Range("D184").Value = 99
$ws.Range("E184").Value = 4
$ws.Range("H184").Value = 3
$ws.Range("B185").Value = 104
$ws.Range("D185").Value = 90
$ws.Range("E185").Value = 7
$ws.Range("H185").Value = 7

# ---- Timestamp update (row 1) ----
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 02:05"
